# InitProperty.xlsx: "unify the conception of DataNode, DataTable, Entity."
#
# The sheet that used to model a single "Property" record is renamed to
# the more general "DataNode", and the cursor/selection is left where the
# author was last working (row 50, column E) instead of the original
# A9 (top of the frozen/scrollable data region).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename worksheet: Property1 -> DataNode
$ws.Name = "DataNode"

# Restore/activate the sheet and move the selection to E50, matching
# where editing left off (header rows 1-8 stay frozen via the existing
# split pane).
$ws.Activate() | Out-Null
$ws.Range("E50").Select() | Out-Null
